$d = $word.ActiveDocument

$old = 'Segundo o texto de Nery et al (2010, p. 2), pl' + [char]0x00E1 + 'gio ' + [char]0x00E9 + ' "copiar" sem citar a fonte, ou n' + [char]0x00E3 + 'o faze-lo corretamente. De acordo com o professor L' + [char]0x00E9 + 'cio Ramos, tamb' + [char]0x00E9 + 'm citado no texto, h' + [char]0x00E1 + ' pelo menos 3 classifica' + [char]0x00E7 + [char]0x00F5 + 'es para pl' + [char]0x00E1 + 'gio, s' + [char]0x00E3 + 'o eles, assim como no texto: pl' + [char]0x00E1 + 'gio integral, quando se copia todo o trabalho, o parcial, que ' + [char]0x00E9 + ' a c' + [char]0x00F3 + 'pia de paragrafos e frases de v' + [char]0x00E1 + 'rios autores e pl' + [char]0x00E1 + 'gio conceitual, ocorre ao utilizar a ideia do autor, porem escrever de maneira diferente, em todos os tipos de pl' + [char]0x00E1 + 'gio a fonte n' + [char]0x00E3 + 'o ' + [char]0x00E9 + ' citada. os autores afirmam que o correto ' + [char]0x00E9 + ' utilizar as pr' + [char]0x00F3 + 'prias palavras, explicar todas as cita' + [char]0x00E7 + [char]0x00F5 + 'es, apresentar as fontes no pr' + [char]0x00F3 + 'prio texto e caso for preciso, fazer cita' + [char]0x00E7 + [char]0x00F5 + 'es diretas.'

$new = 'De acordo com o texto de Nery et al o ato de plagiar ' + [char]0x00E9 + ' copiar ideias, ou conceitos, ou frases, sejam retirados de livros, revistas, ou da internet, de outro autor, o qual as formulou e publicou, e n' + [char]0x00E3 + 'o o citar como fonte de pesquisa. A cartilha aponta o professor L' + [char]0x00E9 + 'cio Ramos, que indica tr' + [char]0x00EA + 's principais ocorr' + [char]0x00EA + 'ncias, s' + [char]0x00E3 + 'o elas, pl' + [char]0x00E1 + 'gio integral o qual configura copiar todo o trabalho, pl' + [char]0x00E1 + 'gio parcial que ' + [char]0x00E9 + ' o agrupamento de trechos e frases de muitos autores, sem citar as fontes, e pl' + [char]0x00E1 + 'gio conceitual que consiste em usar o conceito/ideia de outro autor, mesmo escrevendo diferente do original. Os autores explicam que ' + [char]0x00E9 + ' melhor utilizar palavras pr' + [char]0x00F3 + 'prias em vez de copiar as palavras de outrem, tamb' + [char]0x00E9 + 'm que deve-se explicar todas as cita' + [char]0x00E7 + [char]0x00F5 + 'es utilizadas e apresentar as fontes ao longo do texto.'

# Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,
#          MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)
$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found) {
    throw "Find.Execute failed to locate the target paragraph text"
}
$found
